# Slide 35 ("Where I think Quarto is good (July 2022)") content swap:
#   - Content Placeholder 2 (shape 2) "Blog" bullet becomes the "Python
#     notebook" bullet (new wording for the body line too).
#   - Content Placeholder 3 (shape 3) "Python notebook" bullet becomes the
#     "Blog" bullet (exact previous wording), and its "Presentations" bullet
#     is split into a bold "Presentations: unsure" line plus a new plain
#     "Only if you have" line.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(35)

# --- Content Placeholder 2 -------------------------------------------------
$leftShape = $s.Shapes.Item(2)
$leftText  = $leftShape.TextFrame.TextRange

# Para 3: "Blog: excellent" -> "Python notebook: excellent"
$leftText.Paragraphs(3, 1).Runs(1, 1).Text = "Python notebook: excellent"

# Para 4: blog blurb -> python notebook blurb
$leftText.Paragraphs(4, 1).Runs(1, 1).Text = "Quarto adds to Python notebooks without detracting anything. All you need are a few YAML lines."

# --- Content Placeholder 3 --------------------------------------------------
$rightShape = $s.Shapes.Item(3)
$rightText  = $rightShape.TextFrame.TextRange

# Para 1: "Python notebook: excellent" -> "Blog: excellent"
$rightText.Paragraphs(1, 1).Runs(1, 1).Text = "Blog: excellent"

# Para 2: python notebook blurb -> blog blurb
$rightText.Paragraphs(2, 1).Runs(1, 1).Text = "Quarto allows me to have a scriptable, Python-based blog. I can automate my blog to tweet and post to LinkedIn when I write new articles."

# Para 3: "Presentations: only if you have" -> "Presentations: unsure"
$presPara = $rightText.Paragraphs(3, 1)
$presPara.Runs(1, 1).Text = "Presentations: unsure"

# Insert a new paragraph right after it ("Only if you have"), then strip the
# bold / spacing-before it inherited from the header paragraph it split off
# of so it reads as a normal body line.
$presPara.InsertAfter("`rOnly if you have") | Out-Null
$newPara = $rightText.Paragraphs(4, 1)
$newPara.Font.Bold = 0
$newPara.ParagraphFormat.SpaceBefore = 0
